# Apply the CP05MOAS-GL005 -> CP05MOAS-GL389 rename, fix deployment numbers,
# add engineering data stream refs, and flip the active sheet/selection.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Moorings")
$ws2 = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Moorings sheet ---------------------------------------------------
# Mooring serial number text
$ws1.Range("A2").Value = "CP05MOAS-GL389"
# Deployment number
$ws1.Range("C2").Value = 1

# --- Asset_Cal_Info sheet ----------------------------------------------
# Reference designators (corrected instrument reference designators +
# engineering stream), now tagged with the new glider serial number.
$ws2.Range("A2:A5").Value = "CP05MOAS-GL389-01-ADCPAM000"
$ws2.Range("A7:A10").Value = "CP05MOAS-GL389-02-FLORTM000"
$ws2.Range("A12").Value = "CP05MOAS-GL389-03-CTDGVM000"
$ws2.Range("A14").Value = "CP05MOAS-GL389-04-DOSTAM000"
$ws2.Range("A16").Value = "CP05MOAS-GL389-05-PARADM000"
$ws2.Range("A18").Value = "CP05MOAS-GL389-00-ENG000000"

# Deployment numbers
$ws2.Range("C2:C5").Value = 1
$ws2.Range("C7:C10").Value = 1
$ws2.Range("C12").Value = 1
$ws2.Range("C14").Value = 1
$ws2.Range("C16").Value = 1
$ws2.Range("C18").Value = 1

# --- Active sheet / selections -----------------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("C12").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("C19").Select() | Out-Null
